$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("D2").Value = "28.927.60"

# Row 3
$ws.Range("D3").Value = "1.882.32"
$ws.Range("E3").Value = "  -3.61%  "

# Row 4
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.72%  "

# Row 6
$ws.Range("E6").Value = "  +0.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4585"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.41%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4079"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.64%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07959"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9910"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.93%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.76%  "

# Row 13
$ws.Range("D13").Value = "1.894.82"
$ws.Range("E13").Value = "  -1.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.899"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.83%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.055"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.54%  "

# Row 16
$ws.Range("E16").Value = "  +0.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.42%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06572"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.83%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001025"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.19%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.08%  "

# Row 21
$ws.Range("E21").Value = "  +0.37%  "

# Row 22
$ws.Range("D22").Value = "28.901.70"
$ws.Range("E22").Value = "  -2.92%  "

# Row 23
$ws.Range("E23").Value = "  -4.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "

# Row 25
$ws.Range("E25").Value = "  -3.30%  "

# Row 26
$ws.Range("D26").Value = "2.113.70"
$ws.Range("E26").Value = "  -2.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.42%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.58%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.075"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.15%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.463"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.95%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.018"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09318"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.56%  "

# Row 34
$ws.Range("E34").Value = "  -5.87%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.521"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.275"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.38%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06039"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.74%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02224"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.299"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5773"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.89%  "

# Row 43
$ws.Range("E43").Value = "  -4.89%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07490"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.69%  "

# Row 48
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5441"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.50%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.232"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.895"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.88%  "
